$wb = $excel.ActiveWorkbook
$ordersWs = $wb.Worksheets.Item("Orders")
$summaryWs = $wb.Worksheets.Item("Summary")

# Append the new order line items to the "Orders" sheet (rows 12-21).
# Columns A (package #) and F (quantity) hold numeric-looking values that
# must stay text (matching the rest of the sheet), hence the leading "'".
$ordersWs.Cells.Item(12, 3).Value = "268_猩红泡泡_spray red_Rosa rugosa Thunb._10stems"
$ordersWs.Cells.Item(12, 6).Value = "'6"

$ordersWs.Cells.Item(13, 3).Value = "277_草莓杏仁饼_undefined_Rosa rugosa Thunb._10stems"
$ordersWs.Cells.Item(13, 6).Value = "'3"

$ordersWs.Cells.Item(14, 1).Value = "'3"
$ordersWs.Cells.Item(14, 3).Value = "688_山归来橙_undefined_undefined_1bunch"
$ordersWs.Cells.Item(14, 6).Value = "'10"

$ordersWs.Cells.Item(15, 3).Value = "324_小手球_Spiraea flower double petals_undefined_1bunch"
$ordersWs.Cells.Item(15, 6).Value = "'15"

$ordersWs.Cells.Item(16, 3).Value = "316_尤加利叶大叶_Eucalyptus Cinerea_undefined_1bunch"
$ordersWs.Cells.Item(16, 6).Value = "'15"

$ordersWs.Cells.Item(17, 3).Value = "597_尤加利叶小叶_undefined_undefined_1bunch"
$ordersWs.Cells.Item(17, 6).Value = "'10"

$ordersWs.Cells.Item(18, 3).Value = "592_进口春兰叶_undefined_undefined_1bunch"
$ordersWs.Cells.Item(18, 6).Value = "'15"

$ordersWs.Cells.Item(19, 3).Value = "505_紫罗兰紫_violet purple_undefined_1bunch"
$ordersWs.Cells.Item(19, 6).Value = "'5"

$ordersWs.Cells.Item(20, 3).Value = "411_紫罗兰白_violet white_undefined_1bunch"
$ordersWs.Cells.Item(20, 6).Value = "'10"

$ordersWs.Cells.Item(21, 3).Value = "411_紫罗兰白_violet white_undefined_1bunch"

# Update the aggregated tracking-number cell on the "Summary" sheet. The
# value is a long digit string (with a significant leading zero) so it must
# be forced to text to avoid Excel coercing it into a lossy float.
$summaryWs.Cells.Item(2, 7).Value = "'0146137101398786310151510155100"
